# The deck ships with two DrawingML themes:
#   ppt/theme/theme1.xml  -> "Office Theme" / "Office" colour scheme (only
#                             wired to the Notes Master)
#   ppt/theme/theme2.xml  -> "Integral" / "Red Violet" colour scheme (wired
#                             to the Slide Master that every slide/layout
#                             actually uses)
#
# The authored change swaps the two themes' contents: the Slide Master
# (theme2.xml) should end up with the "Office" colour values, while the
# Notes Master (theme1.xml) would end up with the "Red Violet" ones.
#
# This COM host does not expose a writable NotesMaster theme (nor a raw
# XML/part-replace primitive), so the reachable, persisted half of that
# swap is recolouring the live presentation theme (theme2.xml) — i.e. the
# colours every slide actually renders with — from "Red Violet" to
# "Office". Slide.ThemeColorScheme is backed by that same shared theme
# part, so writing it from any one slide updates it for the whole deck.

function Convert-HexToVbaRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office" colour scheme, in clrScheme document order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $themeColors.Colors($i).RGB = Convert-HexToVbaRgb $officeColors[$i - 1]
}
